$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header label text updates (ms -> s wording)
# ---------------------------------------------------------------
$ws.Range("F3").Value = "BruteForce [s]"
$ws.Range("G3").Value = "Held-Karp [s]"
$ws.Range("P3").Value = "time [s]"
$ws.Range("Q3").Value = "expected time [s]"
$ws.Range("R3").Value = "expected time [s]"

# ---------------------------------------------------------------
# Column F (BruteForce time, column D is raw ms -> column F now in seconds)
# ---------------------------------------------------------------
$ws.Range("F4").Formula = "=D4/10/1000"
$ws.Range("F5:F9").Formula = "=D5/10/1000"
$ws.Range("F4:F8").NumberFormat = "0.0000"
$ws.Range("F9").NumberFormat = "0.000"

# ---------------------------------------------------------------
# Column P (KARP measured time, column N is raw ms -> column P now in seconds)
# ---------------------------------------------------------------
$ws.Range("P7").Formula = "=N7/10/1000"
$ws.Range("P8").Formula = ""
$ws.Range("P9:P20").Formula = "=N9/10/1000"
$ws.Range("P7:P15").NumberFormat = "0.0000"
$ws.Range("P16:P20").NumberFormat = "0.000"
$ws.Range("P17").Formula = ""

# ---------------------------------------------------------------
# Column Q (expected time, theoretical formula) -> seconds
# ---------------------------------------------------------------
$ws.Range("Q7:Q18").Formula = "=2^O7*O7^2/750000/1000"
$ws.Range("Q7:Q10").NumberFormat = "0.0000"

# ---------------------------------------------------------------
# Column O extends down to row 20 (n = 15..28)
# ---------------------------------------------------------------
$ws.Range("O19:O20").Formula = "=O18+1"

# ---------------------------------------------------------------
# Column R (new) - mirrors column Q shifted by two rows, in seconds
# ---------------------------------------------------------------
$ws.Range("R9").Formula = "=Q7"
$ws.Range("R10:R20").Formula = "=Q8"
$ws.Range("R9:R17").NumberFormat = "0.0000"
$ws.Range("R18:R20").NumberFormat = "0.000"

# ---------------------------------------------------------------
# Column G (Held-Karp measured time) references column P, now seconds
# ---------------------------------------------------------------
$ws.Range("G12").Formula = "=P9"
$ws.Range("G16").Formula = "=P9"
$ws.Range("G17").Formula = "=P9"
$ws.Range("G19").Formula = "=P16"
$ws.Range("G21").Formula = "=P16"
$ws.Range("G18").Formula = ""
$ws.Range("G20").Formula = ""
$ws.Range("G12:G17").NumberFormat = "0.0000"
$ws.Range("G18:G21").NumberFormat = "0.000"

# ---------------------------------------------------------------
# Sheet view tweaks: drop frozen/top-left anchor, move selection
# ---------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B27").Select()

# ---------------------------------------------------------------
# Column width for new column R
# ---------------------------------------------------------------
$ws.Range("R1").ColumnWidth = 15.28515625

$wb.Application.Calculate()
